$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "source"
$ws.Range("C2").Value = "mobiles"

$ws.Range("C2").Select()
